$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.395.83"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "3.808.75"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("D7").Value = "3.807.61"
$ws.Range("E7").Value = "  +1.67%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "4.447.74"
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").Value = "3.801.45"
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("D17").Value = "68.434.07"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "465.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.703"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000150"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("E26").Value = "  -2.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "30.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").Value = "3.761.32"
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.40%  "
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +17.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.301"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.32%  "
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.02%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "146.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "392.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("D51").Value = "2.801.04"
$ws.Range("E51").Value = "  +4.39%  "
